$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.774648427963257
$ws.Range("B1").Value = 3.952715635299683
$ws.Range("C1").Value = 3.242544889450073
$ws.Range("D1").Value = 3.093845129013062
$ws.Range("E1").Value = 1.917772173881531
